$wb = $excel.ActiveWorkbook

# --- Step 1: the existing "总计" sheet (2nd tab) is repurposed into "2022-Q1" ---
$wsQ1 = $wb.Worksheets.Item(2)

# Grab the header/index style (bold font + thin border, currently style used by
# the old "总计" header row / A-column index cells) before we touch anything,
# so it can be stamped onto the brand-new "总计" sheet created below.
$wsQ1.Range("B1").Copy()
$headerStyleRange = $wsQ1.Range("B1")

$wsQ1.Name = "2022-Q1"

# New header row for the fund-holding detail sheet (same shape as the "2021-Q3" sheet).
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"

# D1 already carries the bold/border header style; stamp the same style onto the
# newly-introduced E1:H1 header cells so the whole header row is consistent.
$wsQ1.Range("D1").Copy()
$wsQ1.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Data row for the single fund held this quarter.
$wsQ1.Range("A2").Value = 0
$wsQ1.Range("B2").Value = "'004685"
$wsQ1.Range("B2").ClearFormats()
$wsQ1.Range("C2").Value = "金元顺安元启灵活配置混合"
$wsQ1.Range("D2").Value = "'5.00"
$wsQ1.Range("D2").ClearFormats()
$wsQ1.Range("E2").Value = "'75.79"
$wsQ1.Range("E2").ClearFormats()
$wsQ1.Range("F2").Value = "'1.01"
$wsQ1.Range("F2").ClearFormats()
$wsQ1.Range("G2").Value = "'0.0505"
$wsQ1.Range("G2").ClearFormats()
$wsQ1.Range("H2").Value = 5

# --- Step 2: append a brand-new "总计" sheet after the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTotal = $wb.Worksheets.Add($null, $lastSheet)
$wsTotal.Name = "总计"

# Stamp the bold/border header style (captured above) onto the new totals header
# row and onto the A-column index cells.
$headerStyleRange.Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$wsTotal.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.05

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.04
